# Updates cryptos list values (Price / Volume(1h) columns) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell + new text. Price-column values that look like a
# plain number (e.g. "206.64") would otherwise be auto-converted to a
# numeric cell by Excel, so those are written with a leading apostrophe to
# force text, then the quote-prefix style that introduces is reset back to
# "Normal" so the cell keeps the workbook's original (style-less) look.
$updates = @(
    @{ Cell = 'D2'; Value = '26.894.61'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  +0.12%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.546.54'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -1.07%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.28%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '206.64'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +0.39%  '; ForceText = $false }
    @{ Cell = 'E6'; Value = '  +0.01%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +0.32%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.246'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -0.18%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '21.42'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  -1.44%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.0582'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -0.12%  '; ForceText = $false }
    @{ Cell = 'E11'; Value = '  -0.96%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '1.767.67'; ForceText = $false }
    @{ Cell = 'E12'; Value = '  -0.90%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '1.552.23'; ForceText = $false }
    @{ Cell = 'E13'; Value = '  -0.74%  '; ForceText = $false }
    @{ Cell = 'E14'; Value = '  -0.81%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '0.511'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  -0.35%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '26.901.60'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  +0.21%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '61.38'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  +0.41%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '214.40'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  +0.11%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '0.0₃0684'; ForceText = $false }
    @{ Cell = 'E19'; Value = '  +0.93%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '7.22'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -1.80%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '4.01'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -2.71%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '9.16'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -0.21%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '1.93'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -2.88%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '151.81'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -1.45%  '; ForceText = $false }
    @{ Cell = 'E26'; Value = '  -1.11%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '14.85'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -0.52%  '; ForceText = $false }
    @{ Cell = 'E28'; Value = '  +0.28%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '0.103'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  +0.44%  '; ForceText = $false }
    @{ Cell = 'E30'; Value = '  -0.71%  '; ForceText = $false }
    @{ Cell = 'E31'; Value = '  -0.64%  '; ForceText = $false }
    @{ Cell = 'E32'; Value = '  +2.23%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '1.365.87'; ForceText = $false }
    @{ Cell = 'E33'; Value = '  -2.53%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '2.95'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +1.38%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '1.52'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  +0.30%  '; ForceText = $false }
    @{ Cell = 'E36'; Value = '  +4.18%  '; ForceText = $false }
    @{ Cell = 'E37'; Value = '  +0.41%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '0.0164'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  -0.28%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '0.521'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -0.52%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '0.804'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  -1.11%  '; ForceText = $false }
    @{ Cell = 'E42'; Value = '  +5.96%  '; ForceText = $false }
    @{ Cell = 'E43'; Value = '  -0.93%  '; ForceText = $false }
    @{ Cell = 'E44'; Value = '  +2.23%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '63.47'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +0.66%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '1.72'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -2.10%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '1.682.57'; ForceText = $false }
    @{ Cell = 'E47'; Value = '  -0.79%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '85.56'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -0.47%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '0.0507'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  +0.57%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.0₇0972'; ForceText = $false }
    @{ Cell = 'E50'; Value = '  -1.25%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '0.0948'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  +0.31%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $range.Value = "'" + $u.Value
        $range.Style = 'Normal'
    } else {
        $range.Value = $u.Value
    }
}
